$d = $word.ActiveDocument

# The old sentence is split across several runs:
#   "2018 " (bold italic) + "Ημερομηνίες παρατήρησης για τον αστερισμό του "
#   + "Περσεύς" + ": " + "30 Οκτωβρίου-8 Νοεμβρίου και 29 Νοεμβρίου-8 Δεκεμβρίου"
#   (one copy also has a trailing " " run). This exact marker only occurs in
# the four paragraphs that need to be updated, so it's safe to match on.
$marker = "Ημερομηνίες παρατήρησης για τον αστερισμό του"
$newSentence = "Ημερομηνίες παρατήρησης για τον αστερισμό του Hercules: 13-22 Ιουνίου, 12-21 Ιουλίου, 10-19 Αυγούστου"

# Collect the paragraphs that hold the old sentence.
$targets = @()
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains($marker)) {
        $targets += $i
    }
}

# Walk backwards so earlier replacements don't shift the indices of
# paragraphs that still need to be processed.
for ($j = $targets.Count - 1; $j -ge 0; $j--) {
    $idx = $targets[$j]
    $p = $d.Paragraphs($idx)
    $start = $p.Range.Start
    $end = $p.Range.End

    # Exclude the trailing paragraph mark from the range we clear.
    $target = $d.Range($start, $end - 1)
    $target.Delete()

    $ins = $d.Range($start, $start)
    $ins.InsertAfter($newSentence)
}
